# =========================================================================
# Edit script: ThreatCatalogComplete.xlsx
#
# - Renames the "Threats" table's "Precondition" column to "PreCondition"
#   (done by renaming the header cell, which keeps the backing table
#   column definition in sync).
# - Restores/updates the "Threat Components" sheet view (frozen pane /
#   selection).
# - Adds a new "ThreatAgentReplyCategory" worksheet right before the
#   "ThreatAgentAttribute" sheet, colours its tab red, sizes its columns,
#   and fills it with the Id / Reply_id / Category_id lookup table.
# - Adds a sheet-scoped defined name "ThreatAgentReplyCategory" pointing
#   at that table's data range.
# - Leaves the new sheet active/selected last (so it matches the
#   workbook's saved "active tab").
# =========================================================================

$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# 1. Rename the "Threats" table header so the column becomes "PreCondition"
# -------------------------------------------------------------------------
$wsThreatComponents = $wb.Worksheets.Item("Threat Components")
$wsThreatComponents.Range("J1").Value2 = "PreCondition"

# -------------------------------------------------------------------------
# 2. Adjust the sheet view of "Threat Components":
#    - frozen pane restarts scrolled at A2 (top of data)
#    - selection moves to J2 (the column we just renamed)
# -------------------------------------------------------------------------
$wsThreatComponents.Activate()
$excel.ActiveWindow.FreezePanes = $false
$wsThreatComponents.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$wsThreatComponents.Range("J2").Select()

# -------------------------------------------------------------------------
# 3. Insert the new "ThreatAgentReplyCategory" sheet right before
#    "ThreatAgentAttribute".
# -------------------------------------------------------------------------
$wsAttribute = $wb.Worksheets.Item("ThreatAgentAttribute")
$wsReplyCategory = $wb.Worksheets.Add($wsAttribute)
$wsReplyCategory.Name = "ThreatAgentReplyCategory"
$wsReplyCategory.Tab.Color = 255

# Column widths.
$wsReplyCategory.Columns.Item(1).ColumnWidth = 4.1640625
$wsReplyCategory.Columns.Item(2).ColumnWidth = 7.6640625
$wsReplyCategory.Columns.Item(3).ColumnWidth = 19.6640625

# Header row.
$wsReplyCategory.Range("A1").Value2 = "Id"
$wsReplyCategory.Range("B1").Value2 = "Reply_id"
$wsReplyCategory.Range("C1").Value2 = "Category_id"

# Data rows (Id, Reply_id, Category_id) -- 151 rows, A2:C152.
$data = New-Object 'object[,]' 151,3
$data[0,0]=1; $data[0,1]=1; $data[0,2]=1
$data[1,0]=2; $data[1,1]=1; $data[1,2]=2
$data[2,0]=3; $data[2,1]=1; $data[2,2]=3
$data[3,0]=4; $data[3,1]=1; $data[3,2]=4
$data[4,0]=5; $data[4,1]=1; $data[4,2]=5
$data[5,0]=6; $data[5,1]=1; $data[5,2]=6
$data[6,0]=7; $data[6,1]=1; $data[6,2]=7
$data[7,0]=8; $data[7,1]=1; $data[7,2]=8
$data[8,0]=9; $data[8,1]=1; $data[8,2]=9
$data[9,0]=10; $data[9,1]=1; $data[9,2]=10
$data[10,0]=11; $data[10,1]=1; $data[10,2]=11
$data[11,0]=12; $data[11,1]=1; $data[11,2]=12
$data[12,0]=13; $data[12,1]=1; $data[12,2]=13
$data[13,0]=14; $data[13,1]=1; $data[13,2]=14
$data[14,0]=15; $data[14,1]=1; $data[14,2]=15
$data[15,0]=16; $data[15,1]=1; $data[15,2]=16
$data[16,0]=18; $data[16,1]=1; $data[16,2]=18
$data[17,0]=19; $data[17,1]=1; $data[17,2]=19
$data[18,0]=20; $data[18,1]=1; $data[18,2]=20
$data[19,0]=21; $data[19,1]=1; $data[19,2]=21
$data[20,0]=22; $data[20,1]=2; $data[20,2]=1
$data[21,0]=23; $data[21,1]=2; $data[21,2]=2
$data[22,0]=24; $data[22,1]=2; $data[22,2]=3
$data[23,0]=25; $data[23,1]=2; $data[23,2]=13
$data[24,0]=26; $data[24,1]=2; $data[24,2]=16
$data[25,0]=27; $data[25,1]=2; $data[25,2]=20
$data[26,0]=28; $data[26,1]=3; $data[26,2]=13
$data[27,0]=30; $data[27,1]=3; $data[27,2]=20
$data[28,0]=31; $data[28,1]=3; $data[28,2]=4
$data[29,0]=32; $data[29,1]=3; $data[29,2]=5
$data[30,0]=33; $data[30,1]=3; $data[30,2]=6
$data[31,0]=34; $data[31,1]=3; $data[31,2]=8
$data[32,0]=35; $data[32,1]=3; $data[32,2]=10
$data[33,0]=36; $data[33,1]=3; $data[33,2]=14
$data[34,0]=37; $data[34,1]=3; $data[34,2]=15
$data[35,0]=38; $data[35,1]=3; $data[35,2]=16
$data[36,0]=39; $data[36,1]=3; $data[36,2]=18
$data[37,0]=40; $data[37,1]=3; $data[37,2]=7
$data[38,0]=41; $data[38,1]=4; $data[38,2]=1
$data[39,0]=42; $data[39,1]=4; $data[39,2]=2
$data[40,0]=44; $data[40,1]=4; $data[40,2]=3
$data[41,0]=45; $data[41,1]=4; $data[41,2]=4
$data[42,0]=46; $data[42,1]=4; $data[42,2]=5
$data[43,0]=47; $data[43,1]=4; $data[43,2]=6
$data[44,0]=48; $data[44,1]=4; $data[44,2]=7
$data[45,0]=50; $data[45,1]=4; $data[45,2]=8
$data[46,0]=54; $data[46,1]=4; $data[46,2]=9
$data[47,0]=55; $data[47,1]=4; $data[47,2]=10
$data[48,0]=57; $data[48,1]=4; $data[48,2]=11
$data[49,0]=58; $data[49,1]=4; $data[49,2]=12
$data[50,0]=60; $data[50,1]=4; $data[50,2]=13
$data[51,0]=61; $data[51,1]=5; $data[51,2]=13
$data[52,0]=62; $data[52,1]=4; $data[52,2]=16
$data[53,0]=63; $data[53,1]=5; $data[53,2]=20
$data[54,0]=64; $data[54,1]=5; $data[54,2]=5
$data[55,0]=65; $data[55,1]=5; $data[55,2]=6
$data[56,0]=66; $data[56,1]=5; $data[56,2]=8
$data[57,0]=67; $data[57,1]=5; $data[57,2]=11
$data[58,0]=68; $data[58,1]=5; $data[58,2]=12
$data[59,0]=69; $data[59,1]=5; $data[59,2]=14
$data[60,0]=70; $data[60,1]=5; $data[60,2]=21
$data[61,0]=71; $data[61,1]=6; $data[61,2]=13
$data[62,0]=72; $data[62,1]=7; $data[62,2]=16
$data[63,0]=73; $data[63,1]=6; $data[63,2]=20
$data[64,0]=74; $data[64,1]=6; $data[64,2]=7
$data[65,0]=75; $data[65,1]=6; $data[65,2]=10
$data[66,0]=76; $data[66,1]=7; $data[66,2]=13
$data[67,0]=77; $data[67,1]=12; $data[67,2]=16
$data[68,0]=78; $data[68,1]=7; $data[68,2]=20
$data[69,0]=79; $data[69,1]=7; $data[69,2]=1
$data[70,0]=80; $data[70,1]=7; $data[70,2]=2
$data[71,0]=81; $data[71,1]=7; $data[71,2]=3
$data[72,0]=82; $data[72,1]=13; $data[72,2]=16
$data[73,0]=83; $data[73,1]=7; $data[73,2]=8
$data[74,0]=84; $data[74,1]=8; $data[74,2]=13
$data[75,0]=86; $data[75,1]=8; $data[75,2]=20
$data[76,0]=87; $data[76,1]=8; $data[76,2]=4
$data[77,0]=88; $data[77,1]=8; $data[77,2]=10
$data[78,0]=89; $data[78,1]=8; $data[78,2]=9
$data[79,0]=90; $data[79,1]=8; $data[79,2]=18
$data[80,0]=91; $data[80,1]=9; $data[80,2]=13
$data[81,0]=93; $data[81,1]=9; $data[81,2]=20
$data[82,0]=94; $data[82,1]=9; $data[82,2]=9
$data[83,0]=95; $data[83,1]=9; $data[83,2]=10
$data[84,0]=96; $data[84,1]=9; $data[84,2]=15
$data[85,0]=97; $data[85,1]=9; $data[85,2]=18
$data[86,0]=98; $data[86,1]=10; $data[86,2]=13
$data[87,0]=100; $data[87,1]=10; $data[87,2]=20
$data[88,0]=101; $data[88,1]=10; $data[88,2]=14
$data[89,0]=102; $data[89,1]=10; $data[89,2]=18
$data[90,0]=103; $data[90,1]=10; $data[90,2]=15
$data[91,0]=104; $data[91,1]=10; $data[91,2]=19
$data[92,0]=105; $data[92,1]=11; $data[92,2]=13
$data[93,0]=107; $data[93,1]=11; $data[93,2]=20
$data[94,0]=108; $data[94,1]=11; $data[94,2]=12
$data[95,0]=109; $data[95,1]=11; $data[95,2]=15
$data[96,0]=110; $data[96,1]=11; $data[96,2]=19
$data[97,0]=111; $data[97,1]=12; $data[97,2]=13
$data[98,0]=113; $data[98,1]=12; $data[98,2]=20
$data[99,0]=114; $data[99,1]=12; $data[99,2]=1
$data[100,0]=115; $data[100,1]=12; $data[100,2]=2
$data[101,0]=116; $data[101,1]=12; $data[101,2]=3
$data[102,0]=117; $data[102,1]=12; $data[102,2]=4
$data[103,0]=118; $data[103,1]=12; $data[103,2]=9
$data[104,0]=119; $data[104,1]=12; $data[104,2]=10
$data[105,0]=121; $data[105,1]=12; $data[105,2]=18
$data[106,0]=122; $data[106,1]=13; $data[106,2]=13
$data[107,0]=124; $data[107,1]=13; $data[107,2]=20
$data[108,0]=125; $data[108,1]=13; $data[108,2]=1
$data[109,0]=126; $data[109,1]=13; $data[109,2]=2
$data[110,0]=127; $data[110,1]=13; $data[110,2]=3
$data[111,0]=128; $data[111,1]=13; $data[111,2]=5
$data[112,0]=129; $data[112,1]=13; $data[112,2]=9
$data[113,0]=130; $data[113,1]=13; $data[113,2]=10
$data[114,0]=131; $data[114,1]=13; $data[114,2]=14
$data[115,0]=133; $data[115,1]=14; $data[115,2]=13
$data[116,0]=135; $data[116,1]=14; $data[116,2]=20
$data[117,0]=136; $data[117,1]=14; $data[117,2]=6
$data[118,0]=137; $data[118,1]=14; $data[118,2]=7
$data[119,0]=138; $data[119,1]=14; $data[119,2]=8
$data[120,0]=139; $data[120,1]=14; $data[120,2]=11
$data[121,0]=140; $data[121,1]=14; $data[121,2]=12
$data[122,0]=141; $data[122,1]=14; $data[122,2]=21
$data[123,0]=142; $data[123,1]=15; $data[123,2]=13
$data[124,0]=144; $data[124,1]=15; $data[124,2]=20
$data[125,0]=145; $data[125,1]=15; $data[125,2]=6
$data[126,0]=146; $data[126,1]=15; $data[126,2]=7
$data[127,0]=147; $data[127,1]=15; $data[127,2]=8
$data[128,0]=148; $data[128,1]=15; $data[128,2]=11
$data[129,0]=149; $data[129,1]=15; $data[129,2]=21
$data[130,0]=150; $data[130,1]=4; $data[130,2]=14
$data[131,0]=151; $data[131,1]=4; $data[131,2]=15
$data[132,0]=154; $data[132,1]=4; $data[132,2]=18
$data[133,0]=155; $data[133,1]=4; $data[133,2]=19
$data[134,0]=156; $data[134,1]=4; $data[134,2]=20
$data[135,0]=157; $data[135,1]=4; $data[135,2]=21
$data[136,0]=159; $data[136,1]=1; $data[136,2]=17
$data[137,0]=160; $data[137,1]=2; $data[137,2]=17
$data[138,0]=161; $data[138,1]=3; $data[138,2]=17
$data[139,0]=162; $data[139,1]=4; $data[139,2]=17
$data[140,0]=163; $data[140,1]=5; $data[140,2]=17
$data[141,0]=164; $data[141,1]=6; $data[141,2]=17
$data[142,0]=165; $data[142,1]=7; $data[142,2]=17
$data[143,0]=166; $data[143,1]=8; $data[143,2]=17
$data[144,0]=167; $data[144,1]=9; $data[144,2]=17
$data[145,0]=168; $data[145,1]=10; $data[145,2]=17
$data[146,0]=169; $data[146,1]=11; $data[146,2]=17
$data[147,0]=170; $data[147,1]=12; $data[147,2]=17
$data[148,0]=171; $data[148,1]=13; $data[148,2]=17
$data[149,0]=172; $data[149,1]=14; $data[149,2]=17
$data[150,0]=173; $data[150,1]=15; $data[150,2]=17

$wsReplyCategory.Range("A2:C152").Value2 = $data

# -------------------------------------------------------------------------
# 4. Sheet-scoped defined name over the data range.
# -------------------------------------------------------------------------
$wsReplyCategory.Names.Add("ThreatAgentReplyCategory", "=ThreatAgentReplyCategory!`$A`$2:`$C`$152")

# -------------------------------------------------------------------------
# 5. View settings + make this the active sheet/selection last, so the
#    saved workbook view points at it (matches the target "active tab").
# -------------------------------------------------------------------------
$wsReplyCategory.Activate()
$excel.ActiveWindow.Zoom = 200
$wsReplyCategory.Range("D3").Select()
